# Auto-generated edit script: updates FFXIV leve-profit market data values
# across all 8 job tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 3735.6924
$ws.Range("I28").Value = 4175.5713
$ws.Range("J28").Value = 3222.5
$ws.Range("K28").Value = 4175.5713
$ws.Range("L28").Value = 3222.5
$ws.Range("M28").Value = -3690.5713
$ws.Range("N28").Value = -4192.5
# Row 43
$ws.Range("H43").Value = 3196
$ws.Range("J43").Value = 3495
$ws.Range("L43").Value = 3495
$ws.Range("N43").Value = -3633
# Row 74
$ws.Range("H74").Value = 9329.666999999999
$ws.Range("I74").Value = 8999.5
$ws.Range("K74").Value = 8999.5
$ws.Range("M74").Value = -8063.5
# Row 77
$ws.Range("H77").Value = 9329.666999999999
$ws.Range("I77").Value = 8999.5
$ws.Range("K77").Value = 44997.5
$ws.Range("M77").Value = -40317.5
# Row 86
$ws.Range("H86").Value = 1594
$ws.Range("I86").Value = 983
$ws.Range("J86").Value = 1899.5
$ws.Range("K86").Value = 983
$ws.Range("L86").Value = 1899.5
$ws.Range("M86").Value = 140
$ws.Range("N86").Value = -4145.5
# Row 89
$ws.Range("H89").Value = 1594
$ws.Range("I89").Value = 983
$ws.Range("J89").Value = 1899.5
$ws.Range("K89").Value = 4915
$ws.Range("L89").Value = 9497.5
$ws.Range("M89").Value = 701
$ws.Range("N89").Value = -20729.5
# Row 98
$ws.Range("H98").Value = 1155.3334
$ws.Range("I98").Value = 742.5714
$ws.Range("K98").Value = 742.5714
$ws.Range("M98").Value = 755.4286
# Row 100
$ws.Range("H100").Value = 366.55554
$ws.Range("I100").Value = 366.55554
$ws.Range("K100").Value = 366.55554
$ws.Range("M100").Value = 174.44446
# Row 112
$ws.Range("H112").Value = 3005.7144
$ws.Range("J112").Value = 3808.3
$ws.Range("L112").Value = 11424.9
$ws.Range("N112").Value = -13640.9
# Row 113
$ws.Range("H113").Value = 6205.8
$ws.Range("I113").Value = 2895.6
$ws.Range("J113").Value = 9516
$ws.Range("K113").Value = 2895.6
$ws.Range("L113").Value = 9516
$ws.Range("M113").Value = 358.4000000000001
$ws.Range("N113").Value = -16024
# Row 122
$ws.Range("H122").Value = 1155.3334
$ws.Range("I122").Value = 742.5714
$ws.Range("K122").Value = 2227.7142
$ws.Range("M122").Value = 222.2857999999997
# Row 137
$ws.Range("H137").Value = 3570.8125
$ws.Range("I137").Value = 2476.9524
$ws.Range("K137").Value = 7430.8572
$ws.Range("M137").Value = -4880.8572
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 19
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
# Row 32
$ws.Range("H32").Value = 12623
$ws.Range("I32").Value = 4647.8184
$ws.Range("K32").Value = 4647.8184
$ws.Range("M32").Value = -4360.8184
# Row 74
$ws.Range("H74").Value = 2109.3333
$ws.Range("I74").Value = 1164.7273
$ws.Range("K74").Value = 1164.7273
$ws.Range("M74").Value = -290.7273
# Row 77
$ws.Range("H77").Value = 2109.3333
$ws.Range("I77").Value = 1164.7273
$ws.Range("K77").Value = 5823.636500000001
$ws.Range("M77").Value = -1455.636500000001
# Row 97
$ws.Range("H97").Value = 3029.182
$ws.Range("I97").Value = 2047.1428
$ws.Range("K97").Value = 2047.1428
$ws.Range("M97").Value = -1551.1428
# Row 132
$ws.Range("H132").Value = 4788.3335
$ws.Range("I132").Value = 5199.375
$ws.Range("K132").Value = 15598.125
$ws.Range("M132").Value = -13068.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 4479.25
$ws.Range("I94").Value = 3969.3333
$ws.Range("J94").Value = 6009
$ws.Range("K94").Value = 3969.3333
$ws.Range("L94").Value = 6009
$ws.Range("M94").Value = -3518.3333
$ws.Range("N94").Value = -6911
# Row 105
$ws.Range("H105").Value = 3298.4517
$ws.Range("I105").Value = 2556.1333
$ws.Range("J105").Value = 3994.375
$ws.Range("K105").Value = 2556.1333
$ws.Range("L105").Value = 3994.375
$ws.Range("M105").Value = -809.1333
$ws.Range("N105").Value = -7488.375
# Row 134
$ws.Range("H134").Value = 2643.6667
$ws.Range("I134").Value = 2422.2727
$ws.Range("J134").Value = 3617.8
$ws.Range("K134").Value = 7266.8181
$ws.Range("L134").Value = 10853.4
$ws.Range("M134").Value = -4731.8181
$ws.Range("N134").Value = -15923.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6454.0425
$ws.Range("I31").Value = 4255.1055
$ws.Range("J31").Value = 7946.1787
$ws.Range("K31").Value = 4255.1055
$ws.Range("L31").Value = 7946.1787
$ws.Range("M31").Value = -3960.1055
$ws.Range("N31").Value = -8536.1787
# Row 34
$ws.Range("H34").Value = 6454.0425
$ws.Range("I34").Value = 4255.1055
$ws.Range("J34").Value = 7946.1787
$ws.Range("K34").Value = 4255.1055
$ws.Range("L34").Value = 7946.1787
$ws.Range("M34").Value = -4053.1055
$ws.Range("N34").Value = -8350.1787
# Row 58
$ws.Range("H58").Value = 1911.2609
$ws.Range("I58").Value = 1975.409
$ws.Range("K58").Value = 1975.409
$ws.Range("M58").Value = -1772.409
# Row 62
$ws.Range("H62").Value = 19562.264
$ws.Range("J62").Value = 19789.3
$ws.Range("L62").Value = 19789.3
$ws.Range("N62").Value = -21037.3
# Row 65
$ws.Range("H65").Value = 19562.264
$ws.Range("J65").Value = 19789.3
$ws.Range("L65").Value = 98946.5
$ws.Range("N65").Value = -105186.5
# Row 107
$ws.Range("H107").Value = 281.6
$ws.Range("I107").Value = 281.6
$ws.Range("K107").Value = 281.6
$ws.Range("M107").Value = 1638.4
# Row 136
$ws.Range("H136").Value = 1911.2609
$ws.Range("I136").Value = 1975.409
$ws.Range("K136").Value = 5926.227000000001
$ws.Range("M136").Value = -3376.227000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 593.9231
$ws.Range("I12").Value = 640.8333
$ws.Range("J12").Value = 553.7143
$ws.Range("K12").Value = 1922.4999
$ws.Range("L12").Value = 1661.1429
$ws.Range("M12").Value = -1749.4999
$ws.Range("N12").Value = -2007.1429
# Row 107
$ws.Range("H107").Value = 876.3333
$ws.Range("J107").Value = 1009.6
$ws.Range("L107").Value = 3028.8
$ws.Range("N107").Value = -6868.8
# Row 131
$ws.Range("H131").Value = 4223833
$ws.Range("I131").Value = 159665
$ws.Range("J131").Value = 6412231
$ws.Range("K131").Value = 478995
$ws.Range("L131").Value = 19236693
$ws.Range("M131").Value = -473955
$ws.Range("N131").Value = -19246773
# Row 137
$ws.Range("H137").Value = 2618.375
$ws.Range("I137").Value = 2515.6
$ws.Range("K137").Value = 7546.799999999999
$ws.Range("M137").Value = -2446.799999999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 668760.8
$ws.Range("J21").Value = 2006000
$ws.Range("L21").Value = 2006000
$ws.Range("N21").Value = -2006346
# Row 30
$ws.Range("H30").Value = 668760.8
$ws.Range("J30").Value = 2006000
$ws.Range("L30").Value = 2006000
$ws.Range("N30").Value = -2006210
# Row 126
$ws.Range("H126").Value = 3099.2856
$ws.Range("I126").Value = 2489
$ws.Range("J126").Value = 3913
$ws.Range("K126").Value = 7467
$ws.Range("L126").Value = 11739
$ws.Range("M126").Value = -4997
$ws.Range("N126").Value = -16679
# Row 132
$ws.Range("H132").Value = 3651.6428
$ws.Range("I132").Value = 3921.4546
$ws.Range("K132").Value = 11764.3638
$ws.Range("M132").Value = -9234.363799999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 54
$ws.Range("H54").Value = 39500
$ws.Range("J54").Value = 39500
$ws.Range("L54").Value = 39500
$ws.Range("N54").Value = -40788
# Row 82
$ws.Range("H82").Value = 1096.6666
$ws.Range("I82").Value = 947.3333
$ws.Range("J82").Value = 1395.3334
$ws.Range("K82").Value = 947.3333
$ws.Range("L82").Value = 1395.3334
$ws.Range("M82").Value = -586.3333
$ws.Range("N82").Value = -2117.3334
# Row 85
$ws.Range("H85").Value = 1096.6666
$ws.Range("I85").Value = 947.3333
$ws.Range("J85").Value = 1395.3334
$ws.Range("K85").Value = 947.3333
$ws.Range("L85").Value = 1395.3334
$ws.Range("M85").Value = 300.6667
$ws.Range("N85").Value = -3891.3334
# Row 93
$ws.Range("H93").Value = 3917.7144
$ws.Range("J93").Value = 3000
$ws.Range("L93").Value = 3000
$ws.Range("N93").Value = -5496
# Row 122
$ws.Range("H122").Value = 2912.5715
$ws.Range("I122").Value = 2247
$ws.Range("K122").Value = 6741
$ws.Range("M122").Value = -4291
# Row 132
$ws.Range("H132").Value = 2029
$ws.Range("I132").Value = 1907.625
$ws.Range("K132").Value = 5722.875
$ws.Range("M132").Value = -3192.875
# Row 136
$ws.Range("H136").Value = 1687.4445
$ws.Range("I136").Value = 671.2308
$ws.Range("J136").Value = 4329.6
$ws.Range("K136").Value = 2013.6924
$ws.Range("L136").Value = 12988.8
$ws.Range("M136").Value = 536.3075999999999
$ws.Range("N136").Value = -18088.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1047
$ws.Range("I81").Value = 1071.6428
$ws.Range("K81").Value = 2143.2856
$ws.Range("M81").Value = -1082.2856
# Row 84
$ws.Range("H84").Value = 1047
$ws.Range("I84").Value = 1071.6428
$ws.Range("K84").Value = 10716.428
$ws.Range("M84").Value = -5412.428
# Row 132
$ws.Range("H132").Value = 2016.5135
$ws.Range("I132").Value = 1921.1936
$ws.Range("J132").Value = 2509
$ws.Range("K132").Value = 5763.5808
$ws.Range("L132").Value = 7527
$ws.Range("M132").Value = -3233.5808
$ws.Range("N132").Value = -12587
# Row 136
$ws.Range("H136").Value = 791.65
$ws.Range("I136").Value = 709.8889
$ws.Range("K136").Value = 2129.6667
$ws.Range("M136").Value = 420.3332999999998

